$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.334.58"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.864.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.13%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "234.35"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4710"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.43%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2863"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.15%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06568"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.24"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07850"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.10%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "96.83"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.883.34"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6949"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.091"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "267.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.298.03"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.75%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007640"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.37%  "
$ws.Range("E20").Value = "  +0.05%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.113.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.223"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.62%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.153"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.419"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.60"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.86"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.27%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.944"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.361"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09920"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.355"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.458"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.96%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.058"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04729"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.135"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.45%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7006"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.713"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01873"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.32%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.796"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.267"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "72.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.956"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.23%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4172"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.18%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8359"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "103.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "971.44"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.112"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.60%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.114"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "34.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05677"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.49%  "
